$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.947.27"
$ws.Range("E2").Value = "  +2.76%  "

$ws.Range("D3").Value = "3.733.50"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("D7").Value = "3.729.84"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").Value = "4.357.37"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").Value = "3.736.35"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").Value = "68.943.33"
$ws.Range("E17").Value = "  +2.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.80%  "

$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "497.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.82%  "

$ws.Range("E22").Value = "  +4.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000142"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.61%  "

$ws.Range("D34").Value = "3.875.69"
$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("D36").Value = "3.673.57"
$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("E38").Value = "  +1.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.21%  "

$ws.Range("E40").Value = "  -0.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "434.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.66%  "

$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.54%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0352"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.96%  "

$ws.Range("D51").Value = "2.753.62"
$ws.Range("E51").Value = "  -1.25%  "
